# Capstone/data/results/repDiffReport.xlsx - "Add matlab results and poster"
# Updates the repDiffReport table with the new MATLAB-run result values,
# moves the active selection to E10, and restyles the C3:L12 color-scale
# conditional format from a 3-color (min/0/max) scale to a simple
# 2-color (green min -> red max) scale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Updated result values (rows 3-12, columns C-L)
# ---------------------------------------------------------------
$ws.Cells.Item(3, 3).Value = 1678
$ws.Cells.Item(3, 4).Value = 1678
$ws.Cells.Item(3, 5).Value = 1678
$ws.Cells.Item(3, 6).Value = 1680
$ws.Cells.Item(3, 7).Value = 1681
$ws.Cells.Item(3, 8).Value = 1684
$ws.Cells.Item(3, 9).Value = 1686
$ws.Cells.Item(3, 10).Value = 1689
$ws.Cells.Item(3, 11).Value = 1695
$ws.Cells.Item(3, 12).Value = 1707

$ws.Cells.Item(4, 3).Value = 1125
$ws.Cells.Item(4, 4).Value = 1125
$ws.Cells.Item(4, 5).Value = 1128
$ws.Cells.Item(4, 6).Value = 1132
$ws.Cells.Item(4, 7).Value = 1137
$ws.Cells.Item(4, 8).Value = 1150
$ws.Cells.Item(4, 9).Value = 1163
$ws.Cells.Item(4, 10).Value = 1182
$ws.Cells.Item(4, 11).Value = 1209
$ws.Cells.Item(4, 12).Value = 1236

$ws.Cells.Item(5, 3).Value = 541
$ws.Cells.Item(5, 4).Value = 540
$ws.Cells.Item(5, 5).Value = 542
$ws.Cells.Item(5, 6).Value = 552
$ws.Cells.Item(5, 7).Value = 560
$ws.Cells.Item(5, 8).Value = 586
$ws.Cells.Item(5, 9).Value = 608
$ws.Cells.Item(5, 10).Value = 638
$ws.Cells.Item(5, 11).Value = 707
$ws.Cells.Item(5, 12).Value = 785

$ws.Cells.Item(6, 3).Value = 477
$ws.Cells.Item(6, 4).Value = 476
$ws.Cells.Item(6, 5).Value = 479
$ws.Cells.Item(6, 6).Value = 488
$ws.Cells.Item(6, 7).Value = 499
$ws.Cells.Item(6, 8).Value = 529
$ws.Cells.Item(6, 9).Value = 551
$ws.Cells.Item(6, 10).Value = 575
$ws.Cells.Item(6, 11).Value = 643
$ws.Cells.Item(6, 12).Value = 727

$ws.Cells.Item(7, 3).Value = 462
$ws.Cells.Item(7, 4).Value = 463
$ws.Cells.Item(7, 5).Value = 467
$ws.Cells.Item(7, 6).Value = 477
$ws.Cells.Item(7, 7).Value = 488
$ws.Cells.Item(7, 8).Value = 515
$ws.Cells.Item(7, 9).Value = 537
$ws.Cells.Item(7, 10).Value = 561
$ws.Cells.Item(7, 11).Value = 631
$ws.Cells.Item(7, 12).Value = 708

$ws.Cells.Item(8, 3).Value = 455
$ws.Cells.Item(8, 4).Value = 455
$ws.Cells.Item(8, 5).Value = 459
$ws.Cells.Item(8, 6).Value = 471
$ws.Cells.Item(8, 7).Value = 483
$ws.Cells.Item(8, 8).Value = 507
$ws.Cells.Item(8, 9).Value = 531
$ws.Cells.Item(8, 10).Value = 555
$ws.Cells.Item(8, 11).Value = 624
$ws.Cells.Item(8, 12).Value = 702

$ws.Cells.Item(9, 3).Value = 462
$ws.Cells.Item(9, 4).Value = 456
$ws.Cells.Item(9, 5).Value = 458
$ws.Cells.Item(9, 6).Value = 469
$ws.Cells.Item(9, 7).Value = 477
$ws.Cells.Item(9, 8).Value = 501
$ws.Cells.Item(9, 9).Value = 524
$ws.Cells.Item(9, 10).Value = 548
$ws.Cells.Item(9, 11).Value = 617
$ws.Cells.Item(9, 12).Value = 696

$ws.Cells.Item(10, 3).Value = 471
$ws.Cells.Item(10, 4).Value = 451
$ws.Cells.Item(10, 5).Value = 450
$ws.Cells.Item(10, 6).Value = 462
$ws.Cells.Item(10, 7).Value = 465
$ws.Cells.Item(10, 8).Value = 491
$ws.Cells.Item(10, 9).Value = 513
$ws.Cells.Item(10, 10).Value = 539
$ws.Cells.Item(10, 11).Value = 608
$ws.Cells.Item(10, 12).Value = 686

$ws.Cells.Item(11, 3).Value = 539
$ws.Cells.Item(11, 4).Value = 463
$ws.Cells.Item(11, 5).Value = 459
$ws.Cells.Item(11, 6).Value = 462
$ws.Cells.Item(11, 7).Value = 466
$ws.Cells.Item(11, 8).Value = 486
$ws.Cells.Item(11, 9).Value = 508
$ws.Cells.Item(11, 10).Value = 534
$ws.Cells.Item(11, 11).Value = 602
$ws.Cells.Item(11, 12).Value = 678

$ws.Cells.Item(12, 3).Value = 766
$ws.Cells.Item(12, 4).Value = 512
$ws.Cells.Item(12, 5).Value = 464
$ws.Cells.Item(12, 6).Value = 454
$ws.Cells.Item(12, 7).Value = 450
$ws.Cells.Item(12, 8).Value = 475
$ws.Cells.Item(12, 9).Value = 498
$ws.Cells.Item(12, 10).Value = 524
$ws.Cells.Item(12, 11).Value = 594
$ws.Cells.Item(12, 12).Value = 670

# ---------------------------------------------------------------
# Conditional format: replace the 3-stop (min/0/max) color scale on
# C3:L12 with a simple 2-stop green(min) -> red(max) color scale.
# ---------------------------------------------------------------
$cfRange = $ws.Range("C3:L12")
$cfRange.FormatConditions.Delete()
$colorScale = $cfRange.FormatConditions.AddColorScale(2)
# Long color values are BGR-encoded (0x00BBGGRR):
#   RGB(0,176,80)  (green, FF00B050) -> 5287936
#   RGB(192,0,0)   (red,   FFC00000) -> 192
$colorScale.ColorScaleCriteria(1).FormatColor.Color = 5287936
$colorScale.ColorScaleCriteria(2).FormatColor.Color = 192

# ---------------------------------------------------------------
# Move the active selection to E10
# ---------------------------------------------------------------
$ws.Range("E10").Select()
